$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1691842900302115
$ws.Range("C2").Value = 0.595166163141994
$ws.Range("J2").Value = 0.02114803625377644
$ws.Range("P2").Value = 0.1359516616314199
$ws.Range("S2").Value = 0.07854984894259819
$ws.Range("B3").Value = 0.01442307692307692
$ws.Range("C3").Value = 0.04326923076923077
$ws.Range("J3").Value = 0.02884615384615385
$ws.Range("P3").Value = 0.7788461538461539
$ws.Range("S3").Value = 0.1346153846153846
$ws.Range("B6").Value = 0.05194805194805195
$ws.Range("D6").Value = 0.02164502164502164
$ws.Range("E6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.05627705627705628
$ws.Range("J6").Value = 0.2554112554112554
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.1904761904761905
$ws.Range("R6").Value = 0.05194805194805195
$ws.Range("S6").Value = 0.3593073593073593
$ws.Range("B7").Value = 0.1689189189189189
$ws.Range("D7").Value = 0.02027027027027027
$ws.Range("E7").Value = 0.006756756756756757
$ws.Range("F7").Value = 0.1013513513513514
$ws.Range("J7").Value = 0.1081081081081081
$ws.Range("O7").Value = 0.02702702702702703
$ws.Range("Q7").Value = 0.1283783783783784
$ws.Range("R7").Value = 0.04054054054054054
$ws.Range("S7").Value = 0.3986486486486486
$ws.Range("B8").Value = 0.108843537414966
$ws.Range("D8").Value = 0.02267573696145125
$ws.Range("E8").Value = 0.002267573696145125
$ws.Range("F8").Value = 0.07482993197278912
$ws.Range("J8").Value = 0.1020408163265306
$ws.Range("O8").Value = 0.01133786848072562
$ws.Range("Q8").Value = 0.1746031746031746
$ws.Range("R8").Value = 0.08163265306122448
$ws.Range("S8").Value = 0.4217687074829932
$ws.Range("B9").Value = 0.1164658634538153
$ws.Range("D9").Value = 0.0321285140562249
$ws.Range("F9").Value = 0.04417670682730924
$ws.Range("J9").Value = 0.108433734939759
$ws.Range("O9").Value = 0.02409638554216868
$ws.Range("Q9").Value = 0.1967871485943775
$ws.Range("R9").Value = 0.07228915662650602
$ws.Range("S9").Value = 0.4056224899598393
$ws.Range("B10").Value = 0.1104693140794224
$ws.Range("D10").Value = 0.02166064981949458
$ws.Range("F10").Value = 0.05848375451263538
$ws.Range("J10").Value = 0.1220216606498195
$ws.Range("O10").Value = 0.01732851985559567
$ws.Range("Q10").Value = 0.207942238267148
$ws.Range("R10").Value = 0.07148014440433213
$ws.Range("S10").Value = 0.3906137184115523
$ws.Range("G11").Value = 0.1338912133891213
$ws.Range("J11").Value = 0.100418410041841
$ws.Range("K11").Value = 0.196652719665272
$ws.Range("L11").Value = 0.5564853556485355
$ws.Range("S11").Value = 0.01255230125523013
$ws.Range("G12").Value = 0.6811594202898551
$ws.Range("J12").Value = 0.2173913043478261
$ws.Range("K12").Value = 0.01449275362318841
$ws.Range("L12").Value = 0.03623188405797102
$ws.Range("S12").Value = 0.05072463768115942
$ws.Range("G13").Value = 0.6341463414634146
$ws.Range("J13").Value = 0.3170731707317073
$ws.Range("S13").Value = 0.04878048780487805
$ws.Range("F15").Value = 0.01415094339622642
$ws.Range("H15").Value = 0.1415094339622641
$ws.Range("I15").Value = 0.09433962264150944
$ws.Range("J15").Value = 0.3726415094339622
$ws.Range("K15").Value = 0.04716981132075472
$ws.Range("M15").Value = 0.004716981132075472
$ws.Range("O15").Value = 0.07075471698113207
$ws.Range("S15").Value = 0.2547169811320755
$ws.Range("F16").Value = 0.02928870292887029
$ws.Range("H16").Value = 0.1715481171548117
$ws.Range("I16").Value = 0.08786610878661087
$ws.Range("J16").Value = 0.4267782426778243
$ws.Range("K16").Value = 0.09623430962343096
$ws.Range("M16").Value = 0.02510460251046025
$ws.Range("O16").Value = 0.0502092050209205
$ws.Range("S16").Value = 0.1129707112970711
$ws.Range("F17").Value = 0.0148619957537155
$ws.Range("H17").Value = 0.1995753715498939
$ws.Range("I17").Value = 0.1295116772823779
$ws.Range("J17").Value = 0.4331210191082803
$ws.Range("K17").Value = 0.07006369426751592
$ws.Range("M17").Value = 0.01273885350318471
$ws.Range("O17").Value = 0.03609341825902335
$ws.Range("S17").Value = 0.1040339702760085
$ws.Range("F18").Value = 0.01764705882352941
$ws.Range("H18").Value = 0.1882352941176471
$ws.Range("I18").Value = 0.1176470588235294
$ws.Range("J18").Value = 0.4705882352941176
$ws.Range("K18").Value = 0.05294117647058823
$ws.Range("M18").Value = 0.01176470588235294
$ws.Range("O18").Value = 0.05294117647058823
$ws.Range("S18").Value = 0.08823529411764706
$ws.Range("F19").Value = 0.01774691358024691
$ws.Range("H19").Value = 0.1898148148148148
$ws.Range("I19").Value = 0.09876543209876543
$ws.Range("J19").Value = 0.4212962962962963
$ws.Range("K19").Value = 0.08796296296296297
$ws.Range("M19").Value = 0.02160493827160494
$ws.Range("N19").Value = 0.0007716049382716049
$ws.Range("O19").Value = 0.07253086419753087
$ws.Range("S19").Value = 0.08950617283950617
